# Daily attendance processing - 2025-10-23 22:44:17
# Reorders the "Recorded By" (column G) entries for rows whose value is one
# of the known multi-entry strings, reversing the order of the
# comma-separated names/emails (e.g. "System, dnasr281@gmail.com" becomes
# "dnasr281@gmail.com, System").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$map = @{
    "backup@backdoor.com, System, system" = "system, System, backup@backdoor.com";
    "System, dnasr281@gmail.com"          = "dnasr281@gmail.com, System";
    "backup@backdoor.com, System"         = "System, backup@backdoor.com";
}

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = "Recorded By"
    $val = $cell.Value2

    if ($null -ne $val -and $map.ContainsKey($val)) {
        $cell.Value2 = $map[$val]
    }
}
